# Apply weekly data refresh: rows 2-30 (data rows) get their
# Fecha/Variedad/Calidad/Volumen/Precio*/Unidad/Origen/Precio-Kg/Kg-o-Unidades
# values re-shuffled among the existing data rows (columns D and H:Q).
# Columns A,B,C,E,F,G,R are identical across all rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values are copied FROM source row's
# original contents INTO the destination row).
$srcForDst = @{
    2  = 11
    3  = 27
    4  = 17
    5  = 20
    6  = 21
    7  = 6
    8  = 4
    9  = 26
    10 = 3
    11 = 13
    12 = 14
    13 = 18
    14 = 19
    15 = 2
    16 = 25
    17 = 16
    18 = 9
    19 = 30
    20 = 22
    21 = 23
    22 = 24
    23 = 12
    24 = 5
    25 = 10
    26 = 7
    27 = 8
    28 = 28
    29 = 29
    30 = 15
}

# Columns that move together with each logical data row.
$cols = @("D", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q")

# Snapshot the original values of every relevant cell (rows 2-30) before
# writing anything, since several rows are both a source and a destination.
# Note: use .Value2 (not .Value) -- this runtime's COM shim mis-resolves the
# overloaded .Value getter/setter, while .Value2 works correctly for both
# numbers and strings.
$snapshot = @{}
for ($r = 2; $r -le 30; $r++) {
    foreach ($col in $cols) {
        $addr = "$col$r"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Write the shuffled data back out.
foreach ($dst in $srcForDst.Keys) {
    $src = $srcForDst[$dst]
    foreach ($col in $cols) {
        $srcAddr = "$col$src"
        $dstAddr = "$col$dst"
        $ws.Range($dstAddr).Value2 = $snapshot[$srcAddr]
    }
}
